# Blue Ridge Community College Organizations - reformat worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the data held in columns A (Category) and B (Organization Name)
#    for every data row (2-18), since the columns are being swapped:
#    A becomes "Organization Name" and B becomes "Categories".
# ---------------------------------------------------------------------
$swapRows = @{
    2  = @("Arts", "Student Leadership Council")
    3  = @("Arts", "Student Ambassadors Program")
    4  = @("Service", "Blue Ridge PASS Program")
    5  = @("Service", "Social Impact Squad")
    6  = @("Service", "Collegiate FFA (CFFA)")
    7  = @("Professional", "Math Haters Club")
    8  = @("Professional", "Nursing Connections")
    9  = @("Academic", "Phi Theta Kappa (PTK) International Honor Society")
    10 = @("Academic", "STEM Club")
    11 = @("Professional", "Veterinary Technology Club")
    12 = @("Service", "Adventure Club")
    13 = @("Special Interest", "Animanga Club")
    14 = @("Academic", "Blue Ridge Christian Fellowship")
    15 = @("Cultural", "BRCC Diversity Club")
    16 = @("Arts", "Constituting America Club")
    17 = @("General", "Dream, Believe, Achieve (DBA) Club")
    18 = @("Academic", "Prism Club")
}

foreach ($r in $swapRows.Keys) {
    $category = $swapRows[$r][0]
    $orgName  = $swapRows[$r][1]
    $ws.Range("A$r").Value = $orgName
    $ws.Range("B$r").Value = $category
}

# ---------------------------------------------------------------------
# 2) Update the header row for the new column meanings / labels.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Organization Name"
$ws.Range("B1").Value = "Categories"
$ws.Range("C1").Value = "Org URL"
$ws.Range("D1").Value = "Image URL"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Phone"
$ws.Range("H1").Value = "Website"
$ws.Range("I1").Value = "LinkedIn"
$ws.Range("J1").Value = "Instagram"
$ws.Range("K1").Value = "Facebook"
$ws.Range("L1").Value = "Twitter"

# ---------------------------------------------------------------------
# 3) Remove column M (the old "Tiktok Link" column) entirely, shifting
#    nothing else (it was the last column) and shrinking the used
#    range from A1:M18 down to A1:L18.
# ---------------------------------------------------------------------
$ws.Columns.Item(13).Delete()

# ---------------------------------------------------------------------
# 4) Re-apply the column widths to match the new layout: columns A and
#    B swap widths (50 / 18), and the former H:L widths (15,16,15,14,14)
#    are replaced with the narrower 9,10,11,10,9.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 49.17
$ws.Columns.Item(2).ColumnWidth = 17.17
$ws.Columns.Item(8).ColumnWidth = 8.17
$ws.Columns.Item(9).ColumnWidth = 9.17
$ws.Columns.Item(10).ColumnWidth = 10.17
$ws.Columns.Item(11).ColumnWidth = 9.17
$ws.Columns.Item(12).ColumnWidth = 8.17
